$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "trimming_*" flag columns to their "flag_trimming_*" equivalents
$ws.Range("A44").Value = "flag_trimming_hex7_0"
$ws.Range("A45").Value = "flag_trimming_hex7_2"
$ws.Range("A46").Value = "flag_trimming_hex8_0"
$ws.Range("A47").Value = "flag_trimming_hex8_2"
$ws.Range("A48").Value = "flag_trimming_2"

# Update the active selection to match the saved view state
$ws.Range("A49").Select()
